$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 9
$cols = 20
$data = New-Object 'object[,]' $rows,$cols

$data[0,0] = "ECs"
$data[0,1] = "Fn1"
$data[0,2] = "Itgb6"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 27.03890566666666
$data[0,7] = 81.11671699999999
$data[0,8] = 0.07096188219033728
$data[0,9] = 0.07096188219033729
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.128124
$data[0,13] = 0.384372
$data[0,14] = 0.3522399658364659
$data[0,15] = 0.352239965836466
$data[0,16] = 3.464332749636
$data[0,17] = 31.178994746724
$data[0,18] = 0.02499561095841572
$data[0,19] = 0.02499561095841573
$data[1,0] = "ECs"
$data[1,1] = "Fn1"
$data[1,2] = "Itgb6"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 27.03890566666666
$data[1,7] = 81.11671699999999
$data[1,8] = 0.07096188219033728
$data[1,9] = 0.07096188219033729
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.2109236666666666
$data[1,13] = 0.632771
$data[1,14] = 0.5798737562109268
$data[1,15] = 0.5798737562109268
$data[1,16] = 5.703145125867444
$data[1,17] = 51.32830613280699
$data[1,18] = 0.04114893317350815
$data[1,19] = 0.04114893317350816
$data[2,0] = "ECs"
$data[2,1] = "Fn1"
$data[2,2] = "Itgb6"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 27.03890566666666
$data[2,7] = 81.11671699999999
$data[2,8] = 0.07096188219033728
$data[2,9] = 0.07096188219033729
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.024693
$data[2,13] = 0.07407900000000001
$data[2,14] = 0.06788627795260727
$data[2,15] = 0.06788627795260727
$data[2,16] = 0.6676716976270001
$data[2,17] = 6.009045278643
$data[2,18] = 0.004817338058413408
$data[2,19] = 0.004817338058413409
$data[3,0] = "FAPs"
$data[3,1] = "Fn1"
$data[3,2] = "Itgb6"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 345.566579
$data[3,7] = 1036.699737
$data[3,8] = 0.9069174311350353
$data[3,9] = 0.9069174311350354
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.128124
$data[3,13] = 0.384372
$data[3,14] = 0.3522399658364659
$data[3,15] = 0.352239965836466
$data[3,16] = 44.275372367796
$data[3,17] = 398.478351310164
$data[3,18] = 0.3194525649595003
$data[3,19] = 0.3194525649595004
$data[4,0] = "FAPs"
$data[4,1] = "Fn1"
$data[4,2] = "Itgb6"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 345.566579
$data[4,7] = 1036.699737
$data[4,8] = 0.9069174311350353
$data[4,9] = 0.9069174311350354
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.2109236666666666
$data[4,13] = 0.632771
$data[4,14] = 0.5798737562109268
$data[4,15] = 0.5798737562109268
$data[4,16] = 72.88816992013632
$data[4,17] = 655.9935292812269
$data[4,18] = 0.5258976173654375
$data[4,19] = 0.5258976173654376
$data[5,0] = "FAPs"
$data[5,1] = "Fn1"
$data[5,2] = "Itgb6"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 345.566579
$data[5,7] = 1036.699737
$data[5,8] = 0.9069174311350353
$data[5,9] = 0.9069174311350354
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.024693
$data[5,13] = 0.07407900000000001
$data[5,14] = 0.06788627795260727
$data[5,15] = 0.06788627795260727
$data[5,16] = 8.533075535247001
$data[5,17] = 76.797679817223
$data[5,18] = 0.06156724881009756
$data[5,19] = 0.06156724881009757
$data[6,0] = "sCs"
$data[6,1] = "Fn1"
$data[6,2] = "Itgb6"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 8.428738666666668
$data[6,7] = 25.286216
$data[6,8] = 0.0221206866746274
$data[6,9] = 0.02212068667462741
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.128124
$data[6,13] = 0.384372
$data[6,14] = 0.3522399658364659
$data[6,15] = 0.352239965836466
$data[6,16] = 1.079923712928
$data[6,17] = 9.719313416352003
$data[6,18] = 0.007791789918549922
$data[6,19] = 0.007791789918549927
$data[7,0] = "sCs"
$data[7,1] = "Fn1"
$data[7,2] = "Itgb6"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 8.428738666666668
$data[7,7] = 25.286216
$data[7,8] = 0.0221206866746274
$data[7,9] = 0.02212068667462741
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.2109236666666666
$data[7,13] = 0.632771
$data[7,14] = 0.5798737562109268
$data[7,15] = 0.5798737562109268
$data[7,16] = 1.777820464948445
$data[7,17] = 16.000384184536
$data[7,18] = 0.01282720567198119
$data[7,19] = 0.01282720567198119
$data[8,0] = "sCs"
$data[8,1] = "Fn1"
$data[8,2] = "Itgb6"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 8.428738666666668
$data[8,7] = 25.286216
$data[8,8] = 0.0221206866746274
$data[8,9] = 0.02212068667462741
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.024693
$data[8,13] = 0.07407900000000001
$data[8,14] = 0.06788627795260727
$data[8,15] = 0.06788627795260727
$data[8,16] = 0.2081308438960001
$data[8,17] = 1.873177595064
$data[8,18] = 0.001501691084096291
$data[8,19] = 0.001501691084096292

$range = $ws.Range($ws.Cells.Item(2,1), $ws.Cells.Item(1 + $rows, $cols))
$range.Value = $data

Write-Host "Wrote $rows x $cols data to $($range.Address())"
